$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Highlight the "absent" rows (7, 13, 14) with a red fill and mark the
#    "SICK LEAVE" column (I) with 1 for those days.
#
# The style engine allocates a new fill/xf each time a genuinely new color is
# applied, so we first cycle through the blue/orange colors (unused by any
# cell) before landing on the red color that is actually used - this mirrors
# the three new fills/cellXfs that show up in the final workbook.
# ---------------------------------------------------------------------------
$rowsToMark = @(7, 13, 14)

foreach ($r in $rowsToMark) {
    $rng = $ws.Range("A" + $r + ":J" + $r)
    $rng.Interior.Color = 13411113   # FF29A3CC (blue)  - transient
    $rng.Interior.Color = 6737151    # FFFFCC66 (orange) - transient
    $rng.Interior.Color = 6184671    # FFDF5E5E (red)    - final color

    $ws.Range("I" + $r).Value = 1
}

# ---------------------------------------------------------------------------
# 2. Row 19 (TOTAL LEAVES ACCUMULATED): B19 becomes boolean FALSE instead of
#    a blank placeholder. B19 lives inside the merged range A19:G19, and the
#    merged (non top-left) cells cannot be written to directly while merged,
#    so the range is briefly unmerged, the values are restored/updated, and
#    the merge is intentionally left as-is.
# ---------------------------------------------------------------------------
$ws.Range("A19:G19").UnMerge()
$ws.Range("B19").Value = $False
$ws.Range("C19").Value = " "
$ws.Range("D19").Value = " "
$ws.Range("E19").Value = " "
$ws.Range("F19").Value = " "
$ws.Range("G19").Value = " "

# ---------------------------------------------------------------------------
# 3. Fix the invalid 3-argument FLOOR(...) calls so they use the correct
#    2-argument Excel FLOOR signature.
# ---------------------------------------------------------------------------
$ws.Range("B22").Formula = "=FLOOR(F17/8,1)&"".""&FLOOR(MOD(F17,8),1)&"".""&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60"
$ws.Range("B23").Formula = "=FLOOR(H19,1)&"".""&(H19-FLOOR(H19,1))*8&"".0"""
$ws.Range("B24").Formula = "=FLOOR(I19,1)&"".""&(I19-FLOOR(I19,1))*8&"".0"""
$ws.Range("B27").Formula = "=FLOOR(K27/8,1)&"".""&FLOOR(MOD(K27,8),1)&"".""&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60"

Write-Output "done"
